# Edit slide 12 ("Conclusion") content placeholder text per commit diff.
#
# Note: assigning TextRange.Text diffs the new string against the existing
# text (common prefix / common suffix) and can split what was a single run
# into multiple runs whenever the new text shares a prefix or suffix with
# the old text (so the unchanged edges keep their original run identity).
# The target XML keeps every paragraph/run here as a single run, so for
# replacements whose old/new text overlap we first stage a disjoint
# placeholder string (sharing no prefix/suffix with either the old or the
# final text) and only then set the final text - both hops become a full
# single-run replace instead of a partial split, while rPr is preserved.
#
# Also: after any Text assignment that changes a range's length, earlier
# TextRange/Characters() handles into that text may point at stale
# start/length, so re-fetch a fresh Characters() range (at the placeholder's
# own length) before writing the final text into it.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$placeholder = "qzj7pv3mkx9"

# Paragraph 2: "Built a working CIFAR-100 classifier using a CNN and Streamlit."
# -> "I built a working CIFAR-100 classifier using a CNN and Streamlit web app."
# This paragraph has 3 runs: [lead text][Streamlit][trailing "."].
# Update the first and last runs individually (via Characters sub-ranges)
# so the middle "Streamlit" run (which carries err="1") is left untouched.
$para2 = $tr.Paragraphs(2)
$para2Start = $para2.Start
$para2Text = $para2.Text

$oldLead = "Built a working CIFAR-100 classifier using a CNN and "
$newLead = "I built a working CIFAR-100 classifier using a CNN and "
$leadIdx = $para2Text.IndexOf($oldLead)
$leadRange = $tr.Characters($para2Start + $leadIdx, $oldLead.Length)
$leadRange.Text = $placeholder
$leadRange = $tr.Characters($para2Start + $leadIdx, $placeholder.Length)
$leadRange.Text = $newLead

# Re-fetch paragraph 2 text/offsets since the lead run length changed.
$para2 = $tr.Paragraphs(2)
$para2Start = $para2.Start
$para2Text = $para2.Text

$oldTrail = "."
$newTrail = " web app."
$trailIdx = $para2Text.LastIndexOf($oldTrail)
$trailRange = $tr.Characters($para2Start + $trailIdx, $oldTrail.Length)
$trailRange.Text = $placeholder
$trailRange = $tr.Characters($para2Start + $trailIdx, $placeholder.Length)
$trailRange.Text = $newTrail

# Paragraph 3: single run, whole-paragraph replace.
$para3 = $tr.Paragraphs(3)
$para3.Text = $placeholder
$para3 = $tr.Paragraphs(3)
$para3.Text = "I demonstrated how AI/ML models can be integrated into user-friendly applications."

# Paragraph 5: single run, whole-paragraph replace.
$para5 = $tr.Paragraphs(5)
$para5.Text = $placeholder
$para5 = $tr.Paragraphs(5)
$para5.Text = "I will look at enhancing model robustness on current dataset."

# Paragraph 6: single run, whole-paragraph replace.
$para6 = $tr.Paragraphs(6)
$para6.Text = $placeholder
$para6 = $tr.Paragraphs(6)
$para6.Text = "I will aim to expand the app feature for an even better user experience."
